$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell -> new displayed value. The Price/Volume columns are stored as literal
# text in the source workbook (inline strings), not numbers/percentages, so
# each value below is written with a leading apostrophe to force Excel to keep
# it as text even though it looks numeric. The style is then reset to "Normal"
# so no quote-prefix formatting lingers on the cell (exact match to source).
$updates = [ordered]@{
    "D2" = "308.41"
    "E2" = "-1.05%"
    "D3" = "38.27"
    "E3" = "-3.12%"
    "D4" = "5.063"
    "E4" = "-1.79%"
    "D5" = "0.07905"
    "E5" = "-3.41%"
    "D6" = "2.033"
    "E6" = "2.77%"
    "D7" = "4.419"
    "E7" = "4.54%"
    "D8" = "8.247"
    "E8" = "1.26%"
    "D9" = "3.065"
    "E9" = "-8.58%"
    "D10" = "0.9319"
    "E10" = "0.53%"
    "D11" = "0.1278"
    "E11" = "-8.74%"
    "D12" = "0.1886"
    "E12" = "-2.04%"
    "D13" = "0.08775"
    "E13" = "-2.60%"
    "D14" = "0.03435"
    "E14" = "-2.23%"
    "D15" = "0.09752"
    "E15" = "-0.75%"
    "D16" = "0.001407"
    "E16" = "0.74%"
    "D17" = "0.006237"
    "E17" = "3.51%"
    "D18" = "3.562"
    "E18" = "-2.84%"
    "D19" = "0.3453"
    "E19" = "-0.16%"
    "D20" = "0.1288"
    "E20" = "-4.50%"
    "D21" = "5.031"
    "E21" = "8.14%"
    "D22" = "0.2524"
    "E22" = "4.50%"
    "D23" = "0.04339"
    "E23" = "-0.72%"
    "D24" = "0.001221"
    "E24" = "-0.43%"
    "D25" = "0.004625"
    "E25" = "-3.65%"
    "E26" = "176.56%"
    "D39" = "0.02216"
    "E39" = "2.32%"
    "D40" = "0.05031"
    "E40" = "-3.07%"
    "D41" = "0.007577"
    "E41" = "2.44%"
    "D42" = "0.009952"
    "E42" = "1.56%"
    "D43" = "0.1372"
    "E43" = "-0.18%"
    "D44" = "0.002033"
    "E44" = "-3.95%"
    "D45" = "0.008846"
    "E45" = "-10.54%"
    "D46" = "0.00006630"
    "E46" = "3.84%"
    "D47" = "0.00000000754"
    "E47" = "0.66%"
    "D48" = "0.003014"
    "E48" = "9.08%"
    "E49" = "20.80%"
    "D50" = "0.00002112"
    "E50" = "0.66%"
    "D51" = "0.0002011"
    "E51" = "0.66%"
}

$quote = [string][char]39
foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.Value = [string]($quote + $updates[$addr])
    $cell.Style = "Normal"
}

Write-Output ("Applied {0} cell updates" -f $updates.Count)
